$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure these cells keep their original text representation (e.g. trailing
# zeros like "279.00" and literal percent strings like "0.99%") instead of
# being auto-converted to numbers by Excel. Force text format first.

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '279.00'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '0.99%'
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '27.41'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '0.27%'
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '4.836'
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '0.79%'
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '0.06377'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '0.32%'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '7.034'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '1.17%'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '1.284'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '1.00%'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.8924'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '1.74%'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.1525'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '-1.02%'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '14.26%'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07496'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '-0.68%'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.02917'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '-3.54%'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.08987'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '-0.60%'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.001571'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '-0.11%'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.0006398'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '-0.29%'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.006107'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '5.27%'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '3.473'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '0.57%'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '3.301'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '0.03%'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '2.297'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '1.13%'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.1349'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '0.90%'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '3.903'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '-1.09%'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.1504'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '8.96%'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.04398'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '-0.51%'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '0.61%'
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '-1.61%'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.0001650'
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '-14.76%'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.04045'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '-2.65%'
$c = $ws.Range('B41')
$c.NumberFormat = '@'
$c.Value = 'KickToken'
$c = $ws.Range('C41')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.006635'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '-3.12%'
$c = $ws.Range('B42')
$c.NumberFormat = '@'
$c.Value = 'BKEXToken'
$c = $ws.Range('C42')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.1405'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '19.15%'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.002059'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '2.03%'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '-0.28%'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.00005537'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '7.36%'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.561'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '5.01%'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.01846'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '-19.73%'
